$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value ("5.5", "300.0", ...) as literal
# TEXT rather than letting it auto-convert to a number, without leaving a
# residual number-format style on the cell.
function Set-TextNumber($rangeAddr, $value) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# K2 used to hold "Test" - repurpose it as the header of the new LTSD
# Parameters block. The now-unreferenced "Test" shared string is dropped
# automatically when the workbook is saved.
$ws.Range("K2").Value = "LTSD Parameters"

# Row 3: Right / Left column headers
$ws.Range("K3").Value = "Right"
$ws.Range("M3").Value = "Left"

# Row 4: sub headers
$ws.Range("K4").Value = "Threshols"
$ws.Range("L4").Value = "Win"
$ws.Range("M4").Value = "Threshold"
$ws.Range("N4").Value = "Win"

# Row 5: parameter values (Left threshold entered before Right threshold) -
# these look numeric, so they need the text-forcing helper.
Set-TextNumber "M5" "5.0"
Set-TextNumber "K5" "5.5"
Set-TextNumber "L5" "300.0"
Set-TextNumber "N5" "300.0"

# Move the active selection to N5, matching the new focused cell.
$ws.Range("N5").Select()
